$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4498.5713
$ws.Range("I86").Value = 2997.5
$ws.Range("K86").Value = 2997.5
$ws.Range("M86").Value = -1874.5

$ws.Range("H89").Value = 4498.5713
$ws.Range("I89").Value = 2997.5
$ws.Range("K89").Value = 14987.5
$ws.Range("M89").Value = -9371.5

$ws.Range("H96").Value = 1261
$ws.Range("I96").Value = 848
$ws.Range("K96").Value = 2544
$ws.Range("M96").Value = -1171

$ws.Range("H100").Value = 5824.8125
$ws.Range("I100").Value = 5477.778
$ws.Range("J100").Value = 6271
$ws.Range("K100").Value = 5477.778
$ws.Range("L100").Value = 6271
$ws.Range("M100").Value = -4936.778
$ws.Range("N100").Value = -7353

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1911.4375
$ws.Range("I61").Value = 1958.9333
$ws.Range("K61").Value = 1958.9333
$ws.Range("M61").Value = -1746.9333

$ws.Range("H74").Value = 1540.0465
$ws.Range("I74").Value = 1580.7097
$ws.Range("K74").Value = 1580.7097
$ws.Range("M74").Value = -706.7097000000001

$ws.Range("H77").Value = 1540.0465
$ws.Range("I77").Value = 1580.7097
$ws.Range("K77").Value = 7903.548500000001
$ws.Range("M77").Value = -3535.548500000001

$ws.Range("H95").Value = 98634.60000000001
$ws.Range("J95").Value = 98634.60000000001
$ws.Range("L95").Value = 98634.60000000001
$ws.Range("N95").Value = -104126.6

$ws.Range("H102").Value = 3254.9473
$ws.Range("I102").Value = 3254.9473
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 3254.9473
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -1632.9473

$ws.Range("H132").Value = 1732.4651
$ws.Range("I132").Value = 1598.7028
$ws.Range("J132").Value = 2557.3333
$ws.Range("K132").Value = 4796.1084
$ws.Range("L132").Value = 7671.999899999999
$ws.Range("M132").Value = -2266.1084
$ws.Range("N132").Value = -12731.9999

$ws.Range("H136").Value = 1911.4375
$ws.Range("I136").Value = 1958.9333
$ws.Range("K136").Value = 5876.7999
$ws.Range("M136").Value = -3326.7999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 8808.267
$ws.Range("I107").Value = 8047
$ws.Range("K107").Value = 8047
$ws.Range("M107").Value = -6127

$ws.Range("H122").Value = 58083.332
$ws.Range("J122").Value = 58083.332
$ws.Range("L122").Value = 58083.332
$ws.Range("N122").Value = -67883.33199999999

$ws.Range("H134").Value = 2249.3076
$ws.Range("I134").Value = 1788.5264
$ws.Range("K134").Value = 5365.5792
$ws.Range("M134").Value = -2830.5792

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3549.6667
$ws.Range("I16").Value = 3500
$ws.Range("K16").Value = 3500
$ws.Range("M16").Value = -3213

$ws.Range("H31").Value = 1824.1818
$ws.Range("I31").Value = 1789.625
$ws.Range("K31").Value = 1789.625
$ws.Range("M31").Value = -1494.625

$ws.Range("H34").Value = 1824.1818
$ws.Range("I34").Value = 1789.625
$ws.Range("K34").Value = 1789.625
$ws.Range("M34").Value = -1587.625

$ws.Range("H99").Value = 3058.52
$ws.Range("I99").Value = 3344
$ws.Range("K99").Value = 3344
$ws.Range("M99").Value = -1846

$ws.Range("H105").Value = 524.4
$ws.Range("I105").Value = 471.55554
$ws.Range("K105").Value = 471.55554
$ws.Range("M105").Value = 1275.44446

$ws.Range("H113").Value = 3549.6667
$ws.Range("I113").Value = 3500
$ws.Range("K113").Value = 3500
$ws.Range("M113").Value = -1330

$ws.Range("H126").Value = 3058.52
$ws.Range("I126").Value = 3344
$ws.Range("K126").Value = 10032
$ws.Range("M126").Value = -7562

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 7.8
$ws.Range("J12").Value = 9.333333
$ws.Range("L12").Value = 27.999999
$ws.Range("N12").Value = -373.999999

$ws.Range("H113").Value = 2009.1
$ws.Range("J113").Value = 1666
$ws.Range("L113").Value = 4998
$ws.Range("N113").Value = -9338

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3767.2222
$ws.Range("I102").Value = 3629.889
$ws.Range("J102").Value = 4041.889
$ws.Range("K102").Value = 3629.889
$ws.Range("L102").Value = 4041.889
$ws.Range("M102").Value = -2007.889
$ws.Range("N102").Value = -7285.889

$ws.Range("H132").Value = 6465.3423
$ws.Range("I132").Value = 6573.4287
$ws.Range("J132").Value = 6162.7
$ws.Range("K132").Value = 19720.2861
$ws.Range("L132").Value = 18488.1
$ws.Range("M132").Value = -17190.2861
$ws.Range("N132").Value = -23548.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2800.7778
$ws.Range("I40").Value = 2150.2856
$ws.Range("K40").Value = 2150.2856
$ws.Range("M40").Value = -2014.2856

$ws.Range("H41").Value = 37498.332
$ws.Range("J41").Value = 37498.332
$ws.Range("L41").Value = 37498.332
$ws.Range("N41").Value = -38374.332

$ws.Range("H43").Value = 2530000
$ws.Range("J43").Value = 5020000
$ws.Range("L43").Value = 5020000
$ws.Range("N43").Value = -5020386

$ws.Range("H55").Value = 847.7273
$ws.Range("I55").Value = 847.7273
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 847.7273
$ws.Range("L55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -674.7273

$ws.Range("H82").Value = 62502380
$ws.Range("I82").Value = 90911530
$ws.Range("J82").Value = 2261.6
$ws.Range("K82").Value = 90911530
$ws.Range("L82").Value = 2261.6
$ws.Range("M82").Value = -90911169
$ws.Range("N82").Value = -2983.6

$ws.Range("H85").Value = 62502380
$ws.Range("I85").Value = 90911530
$ws.Range("J85").Value = 2261.6
$ws.Range("K85").Value = 90911530
$ws.Range("L85").Value = 2261.6
$ws.Range("M85").Value = -90910282
$ws.Range("N85").Value = -4757.6

$ws.Range("H100").Value = 71101.17999999999
$ws.Range("I100").Value = 90824.62
$ws.Range("K100").Value = 90824.62
$ws.Range("M100").Value = -90283.62

$ws.Range("H109").Value = 21000
$ws.Range("J109").Value = 21000
$ws.Range("L109").Value = 21000
$ws.Range("N109").Value = -23774

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3691
$ws.Range("I62").Value = 1749.6666
$ws.Range("K62").Value = 1749.6666
$ws.Range("M62").Value = -1125.6666

$ws.Range("H65").Value = 3691
$ws.Range("I65").Value = 1749.6666
$ws.Range("K65").Value = 8748.333000000001
$ws.Range("M65").Value = -5628.333000000001

$ws.Range("H100").Value = 1066.6666
$ws.Range("I100").Value = 1690.5
$ws.Range("J100").Value = 567.6
$ws.Range("K100").Value = 3381
$ws.Range("L100").Value = 1135.2
$ws.Range("M100").Value = -2840
$ws.Range("N100").Value = -2217.2

$ws.Range("H105").Value = 31995
$ws.Range("J105").Value = 31995
$ws.Range("L105").Value = 31995
$ws.Range("N105").Value = -38983

$ws.Range("H106").Value = 80000
$ws.Range("J106").Value = 80000
$ws.Range("L106").Value = 80000
$ws.Range("N106").Value = -82524

$ws.Range("H136").Value = 1198.3334
$ws.Range("I136").Value = 797.6
$ws.Range("K136").Value = 2392.8
$ws.Range("M136").Value = 157.1999999999998

$ws.Range("H139").Value = 39750
$ws.Range("J139").Value = 39750
$ws.Range("L139").Value = 39750
$ws.Range("N139").Value = -50030

$ws.Range("H141").Value = 100000
$ws.Range("J141").Value = 100000
$ws.Range("L141").Value = 100000
$ws.Range("N141").Value = -110360
